# Update gh-pages output (generated at 456a3b4)
# Applies the scraped bilibili listing refresh across the four sheets:
#   展览 (Exhibitions), 演出 (Performances), 本地生活 (Local life), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions) - want-to-go counts (+ 2 sold-out/not-for-sale flips)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 394
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F4").Value = 172
$ws1.Range("G4").Value = "已售罄"
$ws1.Range("F5").Value = 1344
$ws1.Range("F7").Value = 2561
$ws1.Range("F8").Value = 944
$ws1.Range("F9").Value = 18919
$ws1.Range("F11").Value = 2001
$ws1.Range("F14").Value = 351
$ws1.Range("F15").Value = 620
$ws1.Range("F17").Value = 214
$ws1.Range("F21").Value = 217
$ws1.Range("F23").Value = 123
$ws1.Range("F24").Value = 6

# ---------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 25

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5916
$ws3.Range("F3").Value = 590

# ---------------------------------------------------------------
# Sheet "全部类型" (All types) - combined/date-sorted feed.
# A new listing ("光与夜之恋 x 线条小狗 x 爱胖达文化" themed restaurant,
# starting 2024-08-17) now sorts ahead of the existing rows 3-5, so those
# three rows cascade down one slot each (their own "want-to-go" / price
# values also refreshed), followed by a batch of like-for-like
# want-to-go-count refreshes further down the sheet.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Column B holds plain "yyyy-mm-dd" text labels (not real dates) in this
# sheet, so force text formatting before assignment to stop Excel's
# autodetect from coercing them into date serials.
$ws4.Range("B3:B5").NumberFormat = "@"

$ws4.Range("B3").Value = "2024-08-17"
$ws4.Range("C3").Value = "广州·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题餐厅"
$ws4.Range("D3").Value = "天河路299号 时尚天河商业广场"
$ws4.Range("E3").Value = "2024.08.17 00:00-10.27 23:59"
$ws4.Range("F3").Value = 5916
$ws4.Range("G3").Value = "已售罄"
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90448"
$ws4.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202408/lwwhTb2q1723430055880.png"

$ws4.Range("B4").Value = "2024-08-19"
$ws4.Range("C4").Value = "广州·排球少年!!垃圾场决战 主题咖啡厅"
$ws4.Range("D4").Value = "多宝街道恩宁路十一甫新街7号 啡约咖啡馆"
$ws4.Range("E4").Value = "2024.08.19 00:00-10.07 23:59"
$ws4.Range("F4").Value = 590
$ws4.Range("G4").Value = 10
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=90613"
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202408/SyeFTEHD1723516066906.png"

$ws4.Range("B5").Value = "2024-08-30"
$ws4.Range("C5").Value = "广州·木灵动漫 二哈和他的白猫师尊主题餐厅"
$ws4.Range("D5").Value = "天河路299号 时尚天河商业广场"
$ws4.Range("E5").Value = "2024.08.30 00:00-09.30 23:59"
$ws4.Range("F5").Value = 564
$ws4.Range("G5").Value = 10
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=91244"
$ws4.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202408/5tZlgklx1724640910069.png"

$ws4.Range("F7").Value = 172
$ws4.Range("G7").Value = "已售罄"
$ws4.Range("F9").Value = 1344
$ws4.Range("F14").Value = 2561
$ws4.Range("F15").Value = 944
$ws4.Range("F16").Value = 18920
$ws4.Range("F22").Value = 2001
$ws4.Range("F25").Value = 351
$ws4.Range("F26").Value = 620
$ws4.Range("F28").Value = 214
$ws4.Range("F35").Value = 217
$ws4.Range("F38").Value = 123
$ws4.Range("F39").Value = 25
$ws4.Range("F41").Value = 6
